# Final prep stuff for submission
#
# Appends a new, empty (single-space) paragraph at the very end of the
# document body -- after the paragraph ending in "...150 animals on the
# screen at once..." and before the section properties (w:sectPr).
# The new paragraph mirrors the run-properties (rFonts cstheme=minorHAnsi)
# of the preceding text but carries no paragraph-level indentation.

$d = $word.ActiveDocument

# Collapsed range sitting exactly at the end of the document's main story,
# i.e. right before the trailing sectPr.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:cstheme="minorHAnsi"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:cstheme="minorHAnsi"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve"> </w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($newParagraphXml)
